# Weekly data refresh: insert two new observation rows at the top of the
# date-ordered block (rows 325-326), pushing the existing rows 325-401 down
# to 327-403. The workbook's used range grows from A1:R401 to A1:R403.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 325 (everything from 325 downward shifts by 2).
$ws.Rows("325:326").Insert()

# New row 325 - same bucket as the old row 325 (Provincia de Chacabuco,
# Primera quality, $5500-$6000 range) but with this week's figures.
$row325 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44508, 13, 100112012, "Espinaca", "Sin especificar", "Primera", 330, 5500, 6000, 5727, "`$/cuna 10 kilos", "Provincia de Chacabuco", 573, 10, "Hortaliza")
for ($i = 0; $i -lt $row325.Length; $i++) {
    $ws.Cells.Item(325, $i + 1).Value = $row325[$i]
}

# New row 326 - same bucket as the old row 326 (Región Metropolitana,
# Primera quality, $5500-$6000 range) but with this week's figures.
$row326 = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44508, 13, 100112012, "Espinaca", "Sin especificar", "Primera", 380, 5500, 6000, 5724, "`$/cuna 10 kilos", "Región Metropolitana", 572, 10, "Hortaliza")
for ($i = 0; $i -lt $row326.Length; $i++) {
    $ws.Cells.Item(326, $i + 1).Value = $row326[$i]
}
